# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5495
$ws1.Range("F4").Value = 350
$ws1.Range("F7").Value = 54
$ws1.Range("F9").Value = 134
$ws1.Range("F10").Value = 328
$ws1.Range("F11").Value = 422
$ws1.Range("F12").Value = 3006
$ws1.Range("F14").Value = 1612

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5495
$ws4.Range("F4").Value = 350
$ws4.Range("F8").Value = 54
$ws4.Range("F10").Value = 134
$ws4.Range("F11").Value = 328
$ws4.Range("F12").Value = 422
$ws4.Range("F13").Value = 3006
$ws4.Range("F15").Value = 1612
